# chore: update Sheets via scheduled runner
# Applies refreshed market-price figures (Universalis snapshot) to the
# Leve profit calculations across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1655.4445
$ws.Range("I6").Value = 1842.7142
$ws.Range("K6").Value = 5528.142599999999
$ws.Range("M6").Value = -5416.142599999999

# Row 15
$ws.Range("H15").Value = 1694.2963
$ws.Range("I15").Value = 1694.2963
$ws.Range("K15").Value = 5082.8889
$ws.Range("M15").Value = -4913.8889

# Row 17
$ws.Range("H17").Value = 3445.7288
$ws.Range("J17").Value = 3462.2456
$ws.Range("L17").Value = 10386.7368
$ws.Range("N17").Value = -10722.7368

# Row 53
$ws.Range("H53").Value = 865.4375
$ws.Range("I53").Value = 428.66666
$ws.Range("J53").Value = 1427
$ws.Range("K53").Value = 428.66666
$ws.Range("L53").Value = 1427
$ws.Range("M53").Value = 208.33334
$ws.Range("N53").Value = -2701

# Row 70
$ws.Range("H70").Value = 25004880
$ws.Range("I70").Value = 2136.3333
$ws.Range("J70").Value = 62508990
$ws.Range("K70").Value = 6408.999899999999
$ws.Range("L70").Value = 187526970
$ws.Range("M70").Value = -6138.999899999999
$ws.Range("N70").Value = -187527510

# Row 73
$ws.Range("H73").Value = 25004880
$ws.Range("I73").Value = 2136.3333
$ws.Range("J73").Value = 62508990
$ws.Range("K73").Value = 6408.999899999999
$ws.Range("L73").Value = 187526970
$ws.Range("M73").Value = -5472.999899999999
$ws.Range("N73").Value = -187528842

# Row 74
$ws.Range("H74").Value = 6889.9
$ws.Range("I74").Value = 6822.1113
$ws.Range("K74").Value = 6822.1113
$ws.Range("M74").Value = -5886.1113

# Row 76
$ws.Range("H76").Value = 7838.7144
$ws.Range("I76").Value = 7988
$ws.Range("K76").Value = 7988
$ws.Range("M76").Value = -7673

# Row 77
$ws.Range("H77").Value = 6889.9
$ws.Range("I77").Value = 6822.1113
$ws.Range("K77").Value = 34110.5565
$ws.Range("M77").Value = -29430.5565

# Row 79
$ws.Range("H79").Value = 7838.7144
$ws.Range("I79").Value = 7988
$ws.Range("K79").Value = 7988
$ws.Range("M79").Value = -6896

# Row 80
$ws.Range("H80").Value = 50009740
$ws.Range("I80").Value = 83334010
$ws.Range("K80").Value = 250002030
$ws.Range("M80").Value = -250001032

# Row 83
$ws.Range("H83").Value = 50009740
$ws.Range("I83").Value = 83334010
$ws.Range("K83").Value = 750006090
$ws.Range("M83").Value = -750001098

# Row 103
$ws.Range("H103").Value = 2178.7144
$ws.Range("J103").Value = 2250.6
$ws.Range("L103").Value = 6751.799999999999
$ws.Range("N103").Value = -7923.799999999999

# Row 132
$ws.Range("H132").Value = 2124.4119
$ws.Range("I132").Value = 2037.303
$ws.Range("K132").Value = 6111.909000000001
$ws.Range("M132").Value = -3581.909000000001

# Row 137
$ws.Range("H137").Value = 1793526.9
$ws.Range("I137").Value = 1162.1666
$ws.Range("K137").Value = 3486.4998
$ws.Range("M137").Value = -936.4998000000001

# Row 138
$ws.Range("H138").Value = 2425.68
$ws.Range("I138").Value = 924.05884
$ws.Range("J138").Value = 2733.241
$ws.Range("K138").Value = 2772.17652
$ws.Range("L138").Value = 8199.723
$ws.Range("M138").Value = 2367.82348
$ws.Range("N138").Value = -18479.723

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19827816
$ws.Range("I32").Value = 22231600
$ws.Range("K32").Value = 22231600
$ws.Range("M32").Value = -22231313

# Row 61
$ws.Range("H61").Value = 2846.9033
$ws.Range("I61").Value = 2407.7
$ws.Range("K61").Value = 2407.7
$ws.Range("M61").Value = -2195.7

# Row 136
$ws.Range("H136").Value = 2846.9033
$ws.Range("I136").Value = 2407.7
$ws.Range("K136").Value = 7223.099999999999
$ws.Range("M136").Value = -4673.099999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 588.17145
$ws.Range("I94").Value = 490.125
$ws.Range("J94").Value = 802.0909
$ws.Range("K94").Value = 490.125
$ws.Range("L94").Value = 802.0909
$ws.Range("M94").Value = -39.125
$ws.Range("N94").Value = -1704.0909

# Row 105
$ws.Range("H105").Value = 3193.7693
$ws.Range("J105").Value = 3321.2727
$ws.Range("L105").Value = 3321.2727
$ws.Range("N105").Value = -6815.2727

# Row 107
$ws.Range("H107").Value = 1208.0769
$ws.Range("J107").Value = 1300.1111
$ws.Range("L107").Value = 1300.1111
$ws.Range("N107").Value = -5140.1111

# Row 134
$ws.Range("H134").Value = 4204915
$ws.Range("I134").Value = 4764604
$ws.Range("K134").Value = 14293812
$ws.Range("M134").Value = -14291277

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 28.173914
$ws.Range("I7").Value = 24
$ws.Range("J7").Value = 56
$ws.Range("K7").Value = 24
$ws.Range("L7").Value = 56
$ws.Range("M7").Value = 89
$ws.Range("N7").Value = -282

# Row 16
$ws.Range("H16").Value = 723.4
$ws.Range("I16").Value = 754.25
$ws.Range("K16").Value = 754.25
$ws.Range("M16").Value = -467.25

# Row 22
$ws.Range("H22").Value = 288.30768
$ws.Range("J22").Value = 499.5
$ws.Range("L22").Value = 499.5
$ws.Range("N22").Value = -1199.5

# Row 58
$ws.Range("H58").Value = 2911.2
$ws.Range("I58").Value = 2525.738
$ws.Range("J58").Value = 4156.5386
$ws.Range("K58").Value = 2525.738
$ws.Range("L58").Value = 4156.5386
$ws.Range("M58").Value = -2322.738
$ws.Range("N58").Value = -4562.5386

# Row 105
$ws.Range("H105").Value = 2485.75
$ws.Range("I105").Value = 1207.25
$ws.Range("K105").Value = 1207.25
$ws.Range("M105").Value = 539.75

# Row 113
$ws.Range("H113").Value = 723.4
$ws.Range("I113").Value = 754.25
$ws.Range("K113").Value = 754.25
$ws.Range("M113").Value = 1415.75

# Row 132
$ws.Range("H132").Value = 4274.857
$ws.Range("I132").Value = 4080.8262
$ws.Range("J132").Value = 4646.75
$ws.Range("K132").Value = 12242.4786
$ws.Range("L132").Value = 13940.25
$ws.Range("M132").Value = -9712.4786
$ws.Range("N132").Value = -19000.25

# Row 134
$ws.Range("H134").Value = 2661.6428
$ws.Range("I134").Value = 2433.0908
$ws.Range("K134").Value = 7299.2724
$ws.Range("M134").Value = -4764.2724

# Row 136
$ws.Range("H136").Value = 2911.2
$ws.Range("I136").Value = 2525.738
$ws.Range("J136").Value = 4156.5386
$ws.Range("K136").Value = 7577.214
$ws.Range("L136").Value = 12469.6158
$ws.Range("M136").Value = -5027.214
$ws.Range("N136").Value = -17569.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 87
$ws.Range("H87").Value = 3000
$ws.Range("I87").Value = 3000
$ws.Range("K87").Value = 9000
$ws.Range("M87").Value = -7752

# Row 90
$ws.Range("H90").Value = 3000
$ws.Range("I90").Value = 3000
$ws.Range("K90").Value = 27000
$ws.Range("M90").Value = -20760

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4264.9
$ws.Range("I61").Value = 4025
$ws.Range("J61").Value = 4624.75
$ws.Range("K61").Value = 4025
$ws.Range("L61").Value = 4624.75
$ws.Range("M61").Value = -3823
$ws.Range("N61").Value = -5028.75

# Row 82
$ws.Range("H82").Value = 2246.8462
$ws.Range("I82").Value = 2032.7894
$ws.Range("J82").Value = 2827.8572
$ws.Range("K82").Value = 2032.7894
$ws.Range("L82").Value = 2827.8572
$ws.Range("M82").Value = -1671.7894
$ws.Range("N82").Value = -3549.8572

# Row 85
$ws.Range("H85").Value = 2246.8462
$ws.Range("I85").Value = 2032.7894
$ws.Range("J85").Value = 2827.8572
$ws.Range("K85").Value = 2032.7894
$ws.Range("L85").Value = 2827.8572
$ws.Range("M85").Value = -784.7893999999999
$ws.Range("N85").Value = -5323.8572

# Row 93
$ws.Range("H93").Value = 3540.3
$ws.Range("I93").Value = 2999.6667
$ws.Range("J93").Value = 3772
$ws.Range("K93").Value = 2999.6667
$ws.Range("L93").Value = 3772
$ws.Range("M93").Value = -1751.6667
$ws.Range("N93").Value = -6268

# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 113
$ws.Range("H113").Value = 4264.9
$ws.Range("I113").Value = 4025
$ws.Range("J113").Value = 4624.75
$ws.Range("K113").Value = 4025
$ws.Range("L113").Value = 4624.75
$ws.Range("M113").Value = -1855
$ws.Range("N113").Value = -8964.75

# Row 122
$ws.Range("H122").Value = 31216.666
$ws.Range("I122").Value = 36421.715
$ws.Range("J122").Value = 12999
$ws.Range("K122").Value = 109265.145
$ws.Range("L122").Value = 38997
$ws.Range("M122").Value = -106815.145
$ws.Range("N122").Value = -43897

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5591.684
$ws.Range("I62").Value = 4648.857
$ws.Range("K62").Value = 4648.857
$ws.Range("M62").Value = -4024.857

# Row 65
$ws.Range("H65").Value = 5591.684
$ws.Range("I65").Value = 4648.857
$ws.Range("K65").Value = 23244.285
$ws.Range("M65").Value = -20124.285

# Row 114
$ws.Range("H114").Value = 34925.668
$ws.Range("J114").Value = 34925.668
$ws.Range("L114").Value = 34925.668
$ws.Range("N114").Value = -43603.668

# Row 126
$ws.Range("H126").Value = 5492.1
$ws.Range("J126").Value = 4001
$ws.Range("L126").Value = 12003
$ws.Range("N126").Value = -16943

# Row 132
$ws.Range("H132").Value = 30988.229
$ws.Range("I132").Value = 34663.87
$ws.Range("J132").Value = 2502
$ws.Range("K132").Value = 103991.61
$ws.Range("L132").Value = 7506
$ws.Range("M132").Value = -101461.61
$ws.Range("N132").Value = -12566

